$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("100EX_Results_Lower_Bounds")

# New label/average rows in column M
$ws.Range("M4").Value = "Av IR:"
$ws.Range("M5").Formula = "=AVERAGE(E2:E17)"

$ws.Range("M7").Value = "Av IR 8:"
$ws.Range("M8").Formula = "=AVERAGE(E10:E17)"

$ws.Range("M10").Value = "Av AP:"
$ws.Range("M11").Formula = "=AVERAGE(H10:H12,H14:H16)"

# Update selection to M12 (also clears the previous topLeftCell scroll position)
$ws.Activate()
[void]$ws.Range("M12").Select()
